$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.763346076011658
$ws.Range("B1").Value = 1.878540635108948
$ws.Range("C1").Value = 1.94161069393158
$ws.Range("D1").Value = 2.525208473205566
$ws.Range("E1").Value = 2.744461059570312
